# Apply the "web elements changed" fix:
#  - Employee sheet, row 4: Peter/Griffin/petgriffin -> Luke/Skywalker/lukesky
#  - NegativeLogins sheet: "Username cannot be empty" / "Password cannot be empty"
#    error messages replaced by a single "Required" message (rows 5-7, col C)
#  - NegativeLogins becomes the active/selected sheet, with new selections

$wb = $excel.ActiveWorkbook

$employee = $wb.Worksheets.Item("Employee")
$negative = $wb.Worksheets.Item("NegativeLogins")

# --- Employee sheet: update row 4 values (was Peter Griffin / petgriffin) ---
$employee.Range("A4").Value = "Luke"
$employee.Range("B4").Value = "Skywalker"
$employee.Range("C4").Value = "lukesky"

# --- NegativeLogins sheet: collapse the two distinct validation messages ---
# into a single "Required" message.
$negative.Range("C5").Value = "Required"
$negative.Range("C6").Value = "Required"
$negative.Range("C7").Value = "Required"

# --- Selections / active sheet ---
$employee.Range("D11").Select() | Out-Null

$negative.Activate()
$negative.Range("C11").Select() | Out-Null
